$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9975323677062988
$ws.Range("B1").Value = 1.042445063591003
$ws.Range("C1").Value = 5.376063346862793
$ws.Range("D1").Value = 1.596341133117676
$ws.Range("E1").Value = 0.9724112749099731
